$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '30.087.18'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -2.87%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.858.94'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -3.85%  '

$ws.Range("E4").Value = '  +0.16%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '233.19'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -3.75%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.12%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4653'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -2.91%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2802'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -2.95%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06524'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -3.99%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.48'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.37%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07812'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.06%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '96.29'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -8.07%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.856.45'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -3.91%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.114'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -3.42%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.6632'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -3.01%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '280.07'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -5.30%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '30.114.61'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -2.51%  '

$ws.Range("E18").Value = '  +0.08%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '5.481'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.84%  '

$ws.Range("E20").Value = '  -2.92%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '2.098.27'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -4.10%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.000007207'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -5.20%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.09%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '6.111'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -4.52%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '9.282'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.94%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '166.47'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.04%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '18.79'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -5.30%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.899'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -10.30%  '

$ws.Range("E29").Value = '  -4.54%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.09520'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -6.23%  '

$ws.Range("E31").Value = '  -4.70%  '

$ws.Range("E32").Value = '  -4.60%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.080'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -6.34%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.04629'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -4.26%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.6985'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -5.56%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.091'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -3.41%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.700'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.91%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01840'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -5.98%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '6.278'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -3.39%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.513'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -4.62%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '72.72'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -5.51%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.8529'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -2.43%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.909'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -6.30%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.9997'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.08%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '103.86'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.26%  '

$ws.Range("E46").Value = '  -5.29%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.003.92'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -2.16%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '7.165'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -5.26%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '9.203'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +1.03%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '33.91'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -3.55%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.1133'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -6.48%  '
